# Quarterly indexing esoteric bug-fix operation
# Update forecast-error summary statistics (ME, MAE, MSE, RMSE, SE, N)
# for quarters Q1..Q9 (rows 2..11) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1335670119460798
$ws.Range("C2").Value = 0.4099671025390884
$ws.Range("D2").Value = 0.3509508125464788
$ws.Range("E2").Value = 0.5924110165640734
$ws.Range("F2").Value = 0.5841535040152176
$ws.Range("G2").Value = 42

$ws.Range("B3").Value = -0.001513401759154079
$ws.Range("C3").Value = 0.475841569153532
$ws.Range("D3").Value = 0.4177735214472426
$ws.Range("E3").Value = 0.6463540217614822
$ws.Range("F3").Value = 0.6543817783518403
$ws.Range("G3").Value = 41

$ws.Range("B4").Value = 0.1309874194782215
$ws.Range("C4").Value = 0.4444290441430878
$ws.Range("D4").Value = 0.3749385748310419
$ws.Range("E4").Value = 0.6123222802014
$ws.Range("F4").Value = 0.6057678922250229
$ws.Range("G4").Value = 40

$ws.Range("B5").Value = 0.02600164070260688
$ws.Range("C5").Value = 0.4864429873837449
$ws.Range("D5").Value = 0.4097178451191367
$ws.Range("E5").Value = 0.6400920598782153
$ws.Range("F5").Value = 0.6479243911420136
$ws.Range("G5").Value = 39

$ws.Range("B6").Value = 0.1413229605261863
$ws.Range("C6").Value = 0.4406483824926739
$ws.Range("D6").Value = 0.3914688192909578
$ws.Range("E6").Value = 0.6256746912661225
$ws.Range("F6").Value = 0.6176868865792924
$ws.Range("G6").Value = 38

$ws.Range("B7").Value = 0.0295561752036937
$ws.Range("C7").Value = 0.4627862106754021
$ws.Range("D7").Value = 0.4228225359431512
$ws.Range("E7").Value = 0.6502480572390441
$ws.Range("F7").Value = 0.6585360833922921
$ws.Range("G7").Value = 37

$ws.Range("B8").Value = 0.1485433660416745
$ws.Range("C8").Value = 0.428250065966257
$ws.Range("D8").Value = 0.4010703637891622
$ws.Range("E8").Value = 0.6333011635779316
$ws.Range("F8").Value = 0.6243668417797383
$ws.Range("G8").Value = 36

$ws.Range("B9").Value = 0.03269697859646408
$ws.Range("C9").Value = 0.4719153317718797
$ws.Range("D9").Value = 0.4377743913932188
$ws.Range("E9").Value = 0.6616452156505168
$ws.Range("F9").Value = 0.6704845803472345
$ws.Range("G9").Value = 35

$ws.Range("B10").Value = 0.08764418001302933
$ws.Range("C10").Value = 0.3538399490260787
$ws.Range("D10").Value = 0.2806086297826358
$ws.Range("E10").Value = 0.5297250511186306
$ws.Range("F10").Value = 0.53028072424651
$ws.Range("G10").Value = 34

$ws.Range("B11").Value = 0.1097826886828203
$ws.Range("C11").Value = 0.441331489103254
$ws.Range("D11").Value = 0.4022870010623976
$ws.Range("E11").Value = 0.6342609881290174
$ws.Range("F11").Value = 0.6343733905601003
$ws.Range("G11").Value = 33
